$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix a single data value: C456 14 -> 34
$ws.Range("C456").Value = 34

# 2. Normalize column B width to match the other data columns (A, C:F),
#    so the column no longer needs its own custom width.
$ws.Columns.Item(2).ColumnWidth = 8

# 3. Append five new data rows (490-494), copying formatting/style from
#    the last existing data row (489) and then filling in the new values.
$newRows = @(
    @{ A = "2024-11-21 01:17:43"; B = "020915.NC"; C = 269; D = 202; E = 32; F = 1 },
    @{ A = "2024-11-21 22:58:51"; B = "020504.NC"; C = 76;  D = 500; E = 6;  F = 1 },
    @{ A = "2024-11-22 01:08:59"; B = "020502.NC"; C = 66;  D = 500; E = 6;  F = 1 },
    @{ A = "2024-11-22 04:53:45"; B = "020500.NC"; C = 71;  D = 500; E = 6;  F = 1 },
    @{ A = "2024-11-22 12:19:13"; B = "020916.NC"; C = 630; D = 199; E = 32; F = 1 }
)

$targetRow = 490
foreach ($rowData in $newRows) {
    $ws.Range("A489:F489").Copy($ws.Range("A" + $targetRow))
    $ws.Rows.Item($targetRow).RowHeight = 13.55

    $ws.Range("A" + $targetRow).Value = $rowData.A
    $ws.Range("B" + $targetRow).Value = $rowData.B
    $ws.Range("C" + $targetRow).Value = $rowData.C
    $ws.Range("D" + $targetRow).Value = $rowData.D
    $ws.Range("E" + $targetRow).Value = $rowData.E
    $ws.Range("F" + $targetRow).Value = $rowData.F

    $targetRow = $targetRow + 1
}
